$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row: "_old" -> "_FV2410", "_new" -> "_FV2504" -----------
$headerRenames = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $headerRenames.Keys) {
    $ws.Range($addr).Value = $headerRenames[$addr]
}

# --- Turn the used range into an Excel Table (ListObject) -----------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- Freeze the header row (split below row 1, panes frozen) --------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
